$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Total" header in column S, row 1
$ws.Range("S1").Value = "Total"

# New totals (sum of B:R) for existing category rows 2-6
$ws.Range("S2").Value = 1944
$ws.Range("S3").Value = 199
$ws.Range("S4").Value = 1075
$ws.Range("S5").Value = 281
$ws.Range("S6").Value = 1253

# New row 7: "Outros" (Others)
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 157
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 15
$ws.Range("E7").Value = 46
$ws.Range("F7").Value = 102
$ws.Range("G7").Value = 93
$ws.Range("H7").Value = 98
$ws.Range("I7").Value = 101
$ws.Range("J7").Value = 106
$ws.Range("K7").Value = 124
$ws.Range("L7").Value = 101
$ws.Range("M7").Value = 108
$ws.Range("N7").Value = 113
$ws.Range("O7").Value = 114
$ws.Range("P7").Value = 116
$ws.Range("Q7").Value = 131
$ws.Range("R7").Value = 429
$ws.Range("S7").Value = 1965

# New row 8: "Total" (Total)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 174
$ws.Range("C8").Value = 16
$ws.Range("D8").Value = 23
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 117
$ws.Range("G8").Value = 113
$ws.Range("H8").Value = 129
$ws.Range("I8").Value = 163
$ws.Range("J8").Value = 224
$ws.Range("K8").Value = 292
$ws.Range("L8").Value = 358
$ws.Range("M8").Value = 471
$ws.Range("N8").Value = 522
$ws.Range("O8").Value = 522
$ws.Range("P8").Value = 692
$ws.Range("Q8").Value = 756
$ws.Range("R8").Value = 2090
$ws.Range("S8").Value = 6717
